# Refresh the "timestamp" column (Z) values for all data rows (2-102).
# These are the re-run timestamps recorded by the pcsmote logging process
# (see commit message: "dataset Us Crime agregado").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z45").Value = "2025-11-13T06:52:45.804469"
$ws.Range("Z46:Z74").Value = "2025-11-13T06:52:46.009853"
$ws.Range("Z75:Z79").Value = "2025-11-13T06:52:46.234212"
$ws.Range("Z80:Z87").Value = "2025-11-13T06:52:46.235213"
$ws.Range("Z88:Z91").Value = "2025-11-13T06:52:46.236218"
$ws.Range("Z92").Value = "2025-11-13T06:52:46.237790"
$ws.Range("Z93:Z95").Value = "2025-11-13T06:52:46.238189"
$ws.Range("Z96:Z101").Value = "2025-11-13T06:52:46.238768"
$ws.Range("Z102").Value = "2025-11-13T06:52:46.239482"
